$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- wipe everything first (values + the old ad-hoc formatting on A2:A5) ---
$ws.Range("A2:G5").Clear()

# --- the new table has two "year" columns (B=2019, D=2018) whose values are
#     otherwise-numeric-looking strings (years, counts, percentages); force
#     those ranges to Text so they land in the file as shared strings
#     (matching the source export) instead of being auto-coerced to numbers. ---
$ws.Range("B1:B8").NumberFormat = "@"
$ws.Range("D1:D8").NumberFormat = "@"

# Row 1 - new header row
$ws.Range("A1").Value = "Unnamed: 0"
$ws.Range("B1").Value = "2019"
$ws.Range("C1").Value = "Unnamed: 1"
$ws.Range("D1").Value = "2018"
$ws.Range("E1").Value = "Unnamed: 2"

# Rows 2-8 - unpivoted data, columns B and D only
$ws.Range("B2").Value = "32,039"
$ws.Range("D2").Value = "32,766"

$ws.Range("B3").Value = "6.9"
$ws.Range("D3").Value = "6.9"

$ws.Range("B4").Value = "38"
$ws.Range("D4").Value = "39"

$ws.Range("B5").Value = "16"
$ws.Range("D5").Value = "16"

$ws.Range("B6").Value = "2"
$ws.Range("D6").Value = "3"

$ws.Range("B7").Value = "40"
$ws.Range("D7").Value = "39"

$ws.Range("B8").Value = "3"
$ws.Range("D8").Value = "3"

# --- re-apply the original bold/bordered/centred header style (still sat on
#     A1 from the source file) across the whole new header row; pasting the
#     format back over B1/D1 also strips the transient "@" number-format
#     applied above so the header row ends up on the very same style index
#     as before, rather than minting a new bold+text combo. ---
$ws.Range("A1").Copy()
$ws.Range("A1:E1").PasteSpecial(-4122)   # xlPasteFormats
$excel.CutCopyMode = $false

# --- the body cells (B2:B8, D2:D8) should carry no explicit style, same as
#     the rest of the original sheet's data cells. ---
$ws.Range("B2:B8").Style = "Normal"
$ws.Range("D2:D8").Style = "Normal"
